$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 19)
$ws.Range("D2").Value = "2021-01-12"
$ws.Range("J2").Value = 160

# Row 4 (was row 21)
$ws.Range("D4").Value = "2022-10-12"
$ws.Range("J4").Value = 250
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 5000
$ws.Range("P4").Value = 312

# Row 5 (was row 6)
$ws.Range("D5").Value = "2022-11-10"
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("P5").Value = 438

# Row 6 (was row 11)
$ws.Range("D6").Value = "2020-12-22"
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("P6").Value = 344

# Row 7 (was row 17)
$ws.Range("D7").Value = "2022-11-18"
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 469

# Row 8 (was row 12)
$ws.Range("D8").Value = "2020-12-23"
$ws.Range("J8").Value = 210
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 344

# Row 9 (was row 20)
$ws.Range("D9").Value = "2022-11-17"
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("P9").Value = 438

# Row 10 (was row 15)
$ws.Range("D10").Value = "2021-02-04"
$ws.Range("J10").Value = 250

# Row 11 (was row 9)
$ws.Range("D11").Value = "2021-02-03"
$ws.Range("J11").Value = 250

# Row 12 (was row 2)
$ws.Range("D12").Value = "2021-02-05"
$ws.Range("J12").Value = 250

# Row 13 (was row 7)
$ws.Range("D13").Value = "2022-10-21"
$ws.Range("J13").Value = 70
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6500
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 406

# Row 14 (was row 16)
$ws.Range("D14").Value = "2020-12-21"
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("P14").Value = 344

# Row 15 (was row 14)
$ws.Range("D15").Value = "2022-11-08"
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("P15").Value = 500

# Row 16 (was row 4)
$ws.Range("D16").Value = "2021-01-08"
$ws.Range("J16").Value = 430

# Row 17 (was row 13)
$ws.Range("D17").Value = "2021-04-06"
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 6000
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 375

# Row 19 (was row 8)
$ws.Range("D19").Value = "2021-02-09"
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 4500
$ws.Range("M19").Value = 4167
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 260

# Row 20 (was row 10)
$ws.Range("D20").Value = "2021-01-14"
$ws.Range("J20").Value = 340
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 5500
$ws.Range("P20").Value = 344

# Row 21 (was row 5)
$ws.Range("D21").Value = "2020-12-24"
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 5500
$ws.Range("P21").Value = 344
